# Add a new task row (row 15) to the Project Plan describing a minor
# update to the SRS document, carried out by Mina Yousry.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the number formats used by sibling rows so the new cells pick up
# the same styles already present in the sheet (date format for the
# Start/Delivery Date columns, non-technical-task format for Task Type).
$ws.Range("D6:E6").Copy()
$ws.Range("D15:E15").PasteSpecial(-4122)

$ws.Range("C13").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("A15").Value = "Mina Yousry"
$ws.Range("B15").Value = "Minor updates to SRS document"
$ws.Range("C15").Value = "Non-Technical "
$ws.Range("D15").Value = (Get-Date -Year 2020 -Month 2 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E15").Value = (Get-Date -Year 2020 -Month 2 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F15").Value = "10 mins"
$ws.Range("G15").Value = "Remove document status and add status table"
$ws.Range("H15").Value = "Pending"

# Update the active selection to match the author's last cursor position.
$ws.Range("H16").Select()
